# [pr4-2907] Remove the unused "Title" column from the Consent sheet/table
# (data + table column + shared-string cleanup), and make the Consent sheet
# the active tab/selection instead of PCNCode.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Consent")

# Delete the entire "Title" column (column A). This shifts ConsentNo,
# Valid from date, Valid to date, Productivity rate - min/max left by one
# column, updates the sheet dimension, and (together with the XML-mapped
# table living on this sheet) drops the now-unused "Title" table column.
$ws.Columns.Item(1).Delete()

# Make "Consent" the active sheet/tab, with the entire first column
# selected (mirrors the post-edit cursor position left behind by deleting
# column A).
$ws.Activate()
$ws.Columns.Item(1).Select()
